# Updated cryptos list on Sun Jun  2 20:35:44 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target text values for each updated cell (prices, % changes, and two
# coin rows whose rank order/data swapped places).
$updates = @{
    'D2' = '67.733.81'
    'E2' = '  -0.03%  '
    'D3' = '3.775.83'
    'E3' = '  -0.72%  '
    'D4' = '1.00'
    'E4' = '  +0.05%  '
    'D5' = '598.82'
    'E5' = '  -0.12%  '
    'D6' = '162.39'
    'E6' = '  -3.07%  '
    'D7' = '3.774.03'
    'E7' = '  -0.80%  '
    'E8' = '  -0.06%  '
    'D9' = '0.512'
    'E9' = '  -1.47%  '
    'D10' = '0.156'
    'E10' = '  -3.23%  '
    'D11' = '0.444'
    'E11' = '  -1.23%  '
    'E12' = '  +4.69%  '
    'D13' = '0.0000244'
    'E13' = '  -3.92%  '
    'D14' = '34.99'
    'E14' = '  -2.72%  '
    'D15' = '4.409.16'
    'D16' = '3.789.22'
    'E16' = '  +1.19%  '
    'D17' = '67.716.77'
    'E17' = '  -0.14%  '
    'D18' = '18.10'
    'E18' = '  -2.07%  '
    'D20' = '6.99'
    'E20' = '  -1.27%  '
    'D21' = '455.55'
    'E21' = '  -1.38%  '
    'E22' = '  -4.88%  '
    'D23' = '0.689'
    'E23' = '  -1.66%  '
    'D24' = '82.84'
    'E24' = '  -0.62%  '
    'D25' = '0.0000142'
    'E25' = '  -6.44%  '
    'D26' = '11.82'
    'E26' = '  -2.24%  '
    'E27' = '  +0.01%  '
    'E28' = '  -1.54%  '
    'D29' = '9.89'
    'E29' = '  -1.27%  '
    'D30' = '3.925.53'
    'E30' = '  -0.64%  '
    'E32' = '  -2.15%  '
    'D33' = '2.58'
    'E33' = '  -7.05%  '
    'D34' = '28.78'
    'E34' = '  -2.50%  '
    'D35' = '1.00'
    'E35' = '  +0.15%  '
    'E36' = '  -1.73%  '
    'D37' = '0.0989'
    'E37' = '  -1.19%  '
    'E38' = '  +4.95%  '
    'D39' = '5.77'
    'E39' = '  -0.16%  '
    'E40' = '  -2.50%  '
    'B41' = 'dogwifhat'
    'C41' = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
    'D41' = '3.15'
    'E41' = '  -8.12%  '
    'B42' = 'FirstDigitalUSD'
    'C42' = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
    'D42' = '1.00'
    'E42' = '  +0.04%  '
    'E43' = '  +0.08%  '
    'D44' = '43.51'
    'E44' = '  +1.35%  '
    'D45' = '47.14'
    'E45' = '  -2.01%  '
    'B46' = 'Monero'
    'C46' = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    'D46' = '151.26'
    'E46' = '  +2.45%  '
    'B47' = 'TheGraph'
    'C47' = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
    'D47' = '0.293'
    'E47' = '  -2.79%  '
    'D48' = '8.26'
    'E48' = '  -0.98%  '
    'D49' = '1.35'
    'E49' = '  -0.92%  '
    'D50' = '1.83'
    'E50' = '  -0.98%  '
    'D51' = '382.67'
    'E51' = '  -2.28%  '
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    # Prefix with an apostrophe so Excel stores the value as literal text
    # (matches the source inlineStr cells) instead of auto-converting
    # number-like strings (e.g. "1.00", "0.689") into numeric values.
    $cell.Value = "'" + $updates[$cellRef]
    # Re-apply the default "Normal" style so the quote-prefix flag does not
    # leave a stray style index different from the cell's original format.
    $cell.Style = "Normal"
}
